$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.995.97'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.829.68'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = '0.9977'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '243.86'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').Value = '0.6317'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('D7').Value = '0.9991'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.07513'
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('D9').Value = '0.2935'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').Value = '22.85'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').Value = '0.07723'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '1.838.15'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '0.6705'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '83.01'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '0.000009713'
$ws.Range('E16').Value = '  +6.68%  '
$ws.Range('D17').Value = '6.081'
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').Value = '29.038.30'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = '12.54'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').Value = '226.27'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '0.9980'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').Value = '7.174'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').Value = '0.9989'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '159.63'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '0.1406'
$ws.Range('E25').Value = '  +3.62%  '
$ws.Range('D26').Value = '8.526'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').Value = '17.89'
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = '1.496'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '4.119'
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').Value = '4.069'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').Value = '1.196'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('D32').Value = '0.05364'
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('D33').Value = '1.857'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').Value = '0.7430'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').Value = '2.651'
$ws.Range('E36').Value = '  +1.80%  '
$ws.Range('D37').Value = '1.244.25'
$ws.Range('E37').Value = '  -3.24%  '
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('D39').Value = '0.01786'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').Value = '6.586'
$ws.Range('E40').Value = '  +3.36%  '
$ws.Range('D41').Value = '0.9014'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').Value = '0.9988'
$ws.Range('D43').Value = '101.55'
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '1.987.76'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('E45').Value = '  +2.88%  '
$ws.Range('D46').Value = '64.71'
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('D47').Value = '0.5100'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').Value = '0.4066'
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('D49').Value = '9.020'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05765'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '6.752'
$ws.Range('E51').Value = '  +1.15%  '
